$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 780 ("「身体の一部として…」" post) entirely, shifting all following rows up by one.
$ws.Rows.Item(780).Delete()
